$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.731.67'
$ws.Range("D3").Value = '1.638.83'
$ws.Range("E3").Value = '  -0.53%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("E6").Value = '  -2.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.26'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.64%  '
$ws.Range("E9").Value = '  +1.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0889'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.06%  '
$ws.Range("D12").Value = '1.870.86'
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("D13").Value = '1.637.46'
$ws.Range("E13").Value = '  -0.55%  '
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.563'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.56%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.76'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("D17").Value = '27.687.56'
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.21'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.88%  '
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("E23").Value = '  +4.77%  '
$ws.Range("E24").Value = '  +3.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.06'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.96'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("E27").Value = '  -0.96%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  +0.45%  '
$ws.Range("E31").Value = '  +0.36%  '
$ws.Range("E32").Value = '  +0.20%  '
$ws.Range("D33").Value = '1.465.52'
$ws.Range("E33").Value = '  +2.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.11'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.96%  '
$ws.Range("E35").Value = '  -2.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.36'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.51%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.569'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.901'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.04'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.99%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("E43").Value = '  -1.64%  '
$ws.Range("E44").Value = '  +1.40%  '
$ws.Range("E45").Value = '  -1.10%  '
$ws.Range("E46").Value = '  -0.59%  '
$ws.Range("D47").Value = '1.780.86'
$ws.Range("E47").Value = '  -0.51%  '
$ws.Range("E48").Value = '  +3.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '87.06'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.35%  '
$ws.Range("E50").Value = '  -1.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0995'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.04%  '
